# Generate Report for handback
# The file "61f88a21-c53c-4b52-91fa-55208aa3f3a2.md" has now been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# Overview sheet's status columns and each locale sheet's Status /
# Latest Handback DateTime columns accordingly.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# Overview sheet: row 3 corresponds to 61f88a21-c53c-4b52-91fa-55208aa3f3a2.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# zh-cn sheet: update Status (B3) and Latest Handback DateTime (G3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $status
$wsZhCn.Range("G3").Value = "2016-01-19 06:53:44"

# de-de sheet: update Status (B3) and Latest Handback DateTime (G3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $status
$wsDeDe.Range("G3").Value = "2016-01-19 06:54:01"
